$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 74 values that changed
$ws.Range("B74").Value = 12004
$ws.Range("F74").Value = 2604
$ws.Range("H74").Value = 2406
$ws.Range("J74").Value = 1126
$ws.Range("K74").Value = 728
$ws.Range("L74").Value = 398
$ws.Range("M74").Value = 6006
$ws.Range("N74").Value = 5790
$ws.Range("Q74").Value = 1718
$ws.Range("R74").Value = 1332
$ws.Range("S74").Value = 10672
$ws.Range("U74").Value = -1407
$ws.Range("V74").Value = -988
$ws.Range("W74").Value = -419
$ws.Range("X74").Value = 2935
$ws.Range("Z74").Value = 2568
$ws.Range("AA74").Value = 108
$ws.Range("AB74").Value = 14
$ws.Range("AC74").Value = 94
$ws.Range("AD74").Value = 8343
$ws.Range("AE74").Value = 7370
$ws.Range("AG74").Value = 965
$ws.Range("AH74").Value = 694

# Add new row 75
# Force the date-like label to be stored as text (matching the other period
# labels in column A), then strip the temporary number format so the cell
# ends up with the default style, like its neighbours.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = 5624
$ws.Range("C75").Value = -13
$ws.Range("D75").Value = -12
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 3971
$ws.Range("G75").Value = -274
$ws.Range("H75").Value = 4245
$ws.Range("I75").Value = -53
$ws.Range("J75").Value = -1431
$ws.Range("K75").Value = 302
$ws.Range("L75").Value = -1733
$ws.Range("M75").Value = 4407
$ws.Range("N75").Value = 4477
$ws.Range("O75").Value = -17
$ws.Range("P75").Value = -53
$ws.Range("Q75").Value = -1257
$ws.Range("R75").Value = 1895
$ws.Range("S75").Value = 3729
$ws.Range("T75").Value = -61
$ws.Range("U75").Value = 3903
$ws.Range("V75").Value = 1918
$ws.Range("W75").Value = 1985
$ws.Range("X75").Value = 84
$ws.Range("Y75").Value = -2879
$ws.Range("Z75").Value = 2963
$ws.Range("AA75").Value = 120
$ws.Range("AB75").Value = 2
$ws.Range("AC75").Value = 119
$ws.Range("AD75").Value = -1313
$ws.Range("AE75").Value = 1905
$ws.Range("AF75").Value = -376
$ws.Range("AG75").Value = -2842
$ws.Range("AH75").Value = 996
